# v10 release edit: Hall counter field-weakening note swap + field-weakening
# offset value bump (33 -> 33.5) on the "Graph" sheet, removal of a couple of
# scratch helper cells, and re-pointing the active tab/selection to the
# "Graph" sheet (from "Ref. Values").

$wb = $excel.ActiveWorkbook

$wsRef   = $wb.Worksheets("Ref. Values")
$wsGraph = $wb.Worksheets("Graph")

# --- Text updates -----------------------------------------------------
# The long explanatory note (previously in A3) and the short "Field
# Weakening offset max" label (previously in D5) trade places; the note's
# numbers are also refreshed (43 -> 44, avg 33 -> 33,5) to match the new
# F7 value below.
$wsGraph.Range("A3").Value = "The offset added to the Hall counter is 23 for the states whit a falling edge of the Hall value and 44 to the state with a rising edge (avg = 33,5)"
$wsGraph.Range("D5").Value = "Field Weakening offset max"

# --- Value update -------------------------------------------------------
# Field weakening offset bumped from 33 to 33.5 counter steps.
$wsGraph.Range("F7").Value = 33.5

# --- Remove now-unused scratch calculations ------------------------------
$wsGraph.Range("I48").ClearContents()
$wsGraph.Range("I49:J49").ClearContents()

# --- Active sheet / selection -------------------------------------------
# The workbook now opens on the "Graph" tab (was "Ref. Values"), with the
# selection left on E7.
$wsGraph.Activate()
$wsGraph.Range("E7").Select()
